# mdl4ui_cg2013.pptx -- "blanck slides for the end of the prez"
#
# 1) Slide 39 ("The plumbing using APT"): merge the two runs
#    "APT is a standard tooling embedded in the JDK " + "since Java 6"
#    into a single run.
# 2) Slide 44 ("Refactoring and Agile practice"): merge the two title
#    runs "Refactoring and " + "Agile practice" into a single run.
# 3) Append a new slide (303) at the end of the deck, same "Title and
#    Content" layout as the other section-title slides, with the title
#    "Not yet presented" and an empty content placeholder.

$p = $ppt.ActivePresentation

# --- 1) slide 39: merge the "APT ... since Java 6" runs -------------------
$s39 = $p.Slides.Item(39)
$body39 = $s39.Shapes.Item(2).TextFrame.TextRange

# paragraph 2 = "APT is a standard tooling embedded in the JDK since Java 6"
#   run 1: chars 93..138  "APT is a standard tooling embedded in the JDK "
#   run 2: chars 139..150 "since Java 6"
$tail39 = $body39.Characters(139, 12)
$tail39.Text = ""
$head39 = $body39.Characters(93, 46)
$head39.Text = "APT is a standard tooling embedded in the JDK since Java 6"

# --- 2) slide 44: merge the "Refactoring and / Agile practice" runs -------
$s44 = $p.Slides.Item(44)
$title44 = $s44.Shapes.Item(1).TextFrame.TextRange

# run 1: chars 1..16  "Refactoring and "
# run 2: chars 17..30 "Agile practice"
$tail44 = $title44.Characters(17, 14)
$tail44.Text = ""
$head44 = $title44.Characters(1, 16)
$head44.Text = "Refactoring and Agile practice"

# --- 3) new slide at the end: "Not yet presented" --------------------------
# Duplicate the (now merged) slide 44 so the new slide inherits the same
# "Title and Content" layout/placeholders, then retitle it and blank the
# content placeholder.
$dup = $s44.Duplicate()
$s45 = $dup.Item(1)

$s45.Shapes.Item(1).TextFrame.TextRange.Text = "Not yet presented"
